# Fixed Partially Matched Crossover.
# Update the "Fitness" column (C) values for generations 0-10 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value  = 4403.692141776349
$ws.Range("C3").Value  = 4403.692141776349
$ws.Range("C4").Value  = 4374.970396024968
$ws.Range("C5").Value  = 4320.451493906685
$ws.Range("C6").Value  = 4003.267882252596
$ws.Range("C7").Value  = 4003.267882252596
$ws.Range("C8").Value  = 4003.267882252596
$ws.Range("C9").Value  = 4003.267882252596
$ws.Range("C10").Value = 4003.267882252596
$ws.Range("C11").Value = 4003.267882252596
$ws.Range("C12").Value = 4003.267882252596
